$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1. Insert the new "Stack:" and "Queue:" sections right after the Linked
#    List's "Delete: O (n)" paragraph, and before its "Search: O (n)" one.
# -------------------------------------------------------------------------

# Locate the "Delete: O (n)" paragraph that belongs to the Linked List
# block (the second occurrence of that exact paragraph text in the doc).
$script:pos = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $txt = $pp.Range.Text
    if ($txt -eq "Delete: O (n)`r") {
        $script:pos = $pp.Range.End
    }
}

function InsertText($text) {
    $rr = $d.Range($script:pos, $script:pos)
    $rr.Text = $text
    $script:pos = $rr.End
}

function CurrentParaIndex() {
    return $d.Range(0, $script:pos).Paragraphs.Count
}

function InsertParagraphBold($text) {
    InsertText ($text + "`r")
    $idx = CurrentParaIndex
    $pp = $d.Paragraphs.Item($idx)
    $pp.Range.Font.Bold = 1
}

# Linked List: Search / Get (new entries, mirroring Array's block)
InsertText "Search: O (n)`r"
InsertText "Get: O ("
InsertText "n"
InsertText ")`r"

# Stack:
InsertParagraphBold "Stack:"
InsertText "Insert: O (1)`r"
InsertText "Delete: O (1"
InsertText ")`r"
InsertText "Search: O (n)`r"
InsertText "Get: O ("
InsertText "n"
InsertText ")`r"

# Queue:
InsertParagraphBold "Queue:"
InsertText "Insert: O (1)`r"
InsertText "Delete: O (1"
InsertText ")"

# -------------------------------------------------------------------------
# 2. Remove the trailing empty paragraphs at the end of the document,
#    keeping only the one that carries the _GoBack bookmark.
# -------------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$bookmarkParaIndex = -1
for ($i = $lastIndex; $i -ge 1; $i--) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -eq "`r") {
        $bookmarkParaIndex = -1
    }
}

# Find the last non-empty paragraph; everything after it (besides the
# bookmark paragraph itself, which is the first empty one following the
# content) should be removed. The bookmark paragraph is the empty
# paragraph immediately after the last piece of real content.
$contentEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -ne "`r") {
        $contentEnd = $i
    }
}

$bookmarkIndex = $contentEnd + 1
$totalParas = $d.Paragraphs.Count
if ($totalParas -gt $bookmarkIndex) {
    $delStart = $d.Paragraphs.Item($bookmarkIndex + 1).Range.Start
    $delEnd = $d.Paragraphs.Item($totalParas).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
